$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "team record" header cells in row 1, styled like the rest of the
# header row (bold font, thin box border, centered / top-aligned).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Every player row (2-47) gets the same team W/L/T record appended.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}
